$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G1: new numeric sample id (3657) for Qiagen Argus X-12 kit column ---
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Cells.Item(1, 7).Value2 = 3657

# --- G4: same allele call as columns E4/F4 ("1,2") ---
$ws.Range("F1").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$ws.Cells.Item(4, 7).Value2 = "1,2"

# --- New loci rows 57-68: fill column A (locus names) first so the shared
#     strings for the locus names are appended before the allele-call strings,
#     matching the order new loci/alleles were introduced. ---
$locusNames = @("DXS7132","DXS7423","DXS8378","DXS10074","DXS10079","DXS10101","DXS10103","DXS10134","DXS10135","DXS10146","DXS10148","HPRTB")
for ($i = 0; $i -lt $locusNames.Length; $i++) {
    $row = 57 + $i
    $ws.Cells.Item($row, 1).Value2 = $locusNames[$i]
}

# Paint column G (style matching B:F columns) for the new rows before setting values.
$ws.Range("F1").Copy()
$ws.Range("G57:G68").PasteSpecial(-4122)

# Column B values (mix of plain numbers, existing shared strings, and brand-new shared strings)
$ws.Cells.Item(57, 2).Value2 = 12
$ws.Cells.Item(58, 2).Value2 = "14,15"
$ws.Cells.Item(59, 2).Value2 = "10,11"
$ws.Cells.Item(60, 2).Value2 = "16,19"
$ws.Cells.Item(61, 2).Value2 = "20,23"
$ws.Cells.Item(62, 2).Value2 = "30,31"
$ws.Cells.Item(63, 2).Value2 = 17
$ws.Cells.Item(64, 2).Value2 = "35,36"
$ws.Cells.Item(65, 2).Value2 = "21.1,27"
$ws.Cells.Item(66, 2).Value2 = 28
$ws.Cells.Item(67, 2).Value2 = "22.1,23.1"
$ws.Cells.Item(68, 2).Value2 = 14

# Column C values (all plain numbers)
$ws.Cells.Item(57, 3).Value2 = 13
$ws.Cells.Item(58, 3).Value2 = 14
$ws.Cells.Item(59, 3).Value2 = 11
$ws.Cells.Item(60, 3).Value2 = 18
$ws.Cells.Item(61, 3).Value2 = 19
$ws.Cells.Item(62, 3).Value2 = 32
$ws.Cells.Item(63, 3).Value2 = 18
$ws.Cells.Item(64, 3).Value2 = 34
$ws.Cells.Item(65, 3).Value2 = 22
$ws.Cells.Item(66, 3).Value2 = 29
$ws.Cells.Item(67, 3).Value2 = 23
$ws.Cells.Item(68, 3).Value2 = 14

# Column D values (all plain numbers)
$ws.Cells.Item(57, 4).Value2 = 13
$ws.Cells.Item(58, 4).Value2 = 17
$ws.Cells.Item(59, 4).Value2 = 10
$ws.Cells.Item(60, 4).Value2 = 17
$ws.Cells.Item(61, 4).Value2 = 17
$ws.Cells.Item(62, 4).Value2 = 31
$ws.Cells.Item(63, 4).Value2 = 17
$ws.Cells.Item(64, 4).Value2 = 32
$ws.Cells.Item(65, 4).Value2 = 27
$ws.Cells.Item(66, 4).Value2 = 29
$ws.Cells.Item(67, 4).Value2 = 23.1
$ws.Cells.Item(68, 4).Value2 = 13

# Column G values (new Qiagen Argus X-12 kit control calls, all plain numbers)
$ws.Cells.Item(57, 7).Value2 = 12
$ws.Cells.Item(58, 7).Value2 = 13
$ws.Cells.Item(59, 7).Value2 = 12
$ws.Cells.Item(60, 7).Value2 = 7
$ws.Cells.Item(61, 7).Value2 = 19
$ws.Cells.Item(62, 7).Value2 = 29.2
$ws.Cells.Item(63, 7).Value2 = 20
$ws.Cells.Item(64, 7).Value2 = 34
$ws.Cells.Item(65, 7).Value2 = 25
$ws.Cells.Item(66, 7).Value2 = 27
$ws.Cells.Item(67, 7).Value2 = 23.1
$ws.Cells.Item(68, 7).Value2 = 13

# --- Adjust the sheet selection to match the saved view state ---
[void]$ws.Range("I63").Select()

Write-Output "done"
